$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Most Probable" / "Least Probable" keep their original Arial/Arail
# labels (A2/A3 formatting is unchanged).

# The letter/number samples are replaced by picture filenames, rendered
# in the plain default font instead of the old oversized custom fonts.
$ws.Range("B2:C3").Font.Bold = $false
$ws.Range("B2:C3").Font.Size = 11
$ws.Range("B2:C3").Font.Name = "Calibri"

$ws.Range("C2").Value = "Nselect2.jpg"
$ws.Range("B2").Value = "Hselect2.jpg"
$ws.Range("C3").Value = "Nselect2.jpg"
$ws.Range("B3").Value = "Hselect2.jpg"

# Rows no longer need to be tall enough for size-60 text in every column.
$ws.Range("A2:A3").RowHeight = 37.5

# Selection / active cell moved.
$ws.Range("B8").Select() | Out-Null

# Page setup now explicit.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
